# ML model retrained with all data.
# Column J (rows 1-51) and cell K1 previously held stale predictions
# (J1/K1 were text "r"/"s" placeholders, J2:J51 were 0.5). The retrained
# model now scores every row at 0.6, so J1:J51 and K1 are rewritten as the
# numeric value 0.6 (K2:K51 already read 0.6 and are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J1 and K1 switch from shared-string labels to the model's numeric output.
$ws.Range("J1:K1").Value = 0.6

# J2:J51 move from the old 0.5 score to the retrained model's 0.6 score.
$ws.Range("J2:J51").Value = 0.6

# Update the view to match where the analyst last left the sheet: scrolled
# down a couple more rows, zoomed to 100%, with K1:K51 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 39
$win.ScrollColumn = 1
$win.Zoom = 100

$ws.Range("K1:K51").Select() | Out-Null
